# Expand the abbreviated pack-code labels (column "code_pack") into their
# full descriptive names. Using whole-cell Replace (not a Value overwrite)
# keeps the shared-string table compact and avoids any accidental partial
# matches (e.g. "DC" is a substring of "SICEDC"/"SPCEDC").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

$ws.Cells.Replace("DG", "Dynamic Gold ", $xlWhole)
$ws.Cells.Replace("DC", "Dynamic Classic", $xlWhole)
$ws.Cells.Replace("SICCH", "Signature Infinite CCH", $xlWhole)
$ws.Cells.Replace("DT", "Dynamic Titanium ", $xlWhole)
$ws.Cells.Replace("SPCCH", "Signature Platinium CCH", $xlWhole)
$ws.Cells.Replace("SICEDC", "Signature Infinite CEDC", $xlWhole)
$ws.Cells.Replace("SG", "Signature Gold", $xlWhole)
$ws.Cells.Replace("ST", "Signature Titanium", $xlWhole)
$ws.Cells.Replace("SPCEDC", "Signature Platinium CEDC", $xlWhole)
